$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 288 - this shifts the existing rows 288..328 down to
# 289..329, preserving all of their data/formatting intact (matches the
# diff's observed "row N (before) -> row N+1 (after)" shift for rows
# 288..327, with the former row 328 becoming row 329).
$ws.Rows.Item(288).Insert()

# Populate the newly inserted row 288 with the new weekly price record.
$ws.Range("A288").Value = 3
$ws.Range("B288").Value = "Femacal de La Calera"
$ws.Range("C288").Value = "Coquimbo"
$ws.Range("D288").Value = 44776
$ws.Range("E288").Value = 5
$ws.Range("F288").Value = 100112001
$ws.Range("G288").Value = "Berenjena"
$ws.Range("H288").Value = "Sin especificar"
$ws.Range("I288").Value = "Primera"
$ws.Range("J288").Value = 150
$ws.Range("K288").Value = 9000
$ws.Range("L288").Value = 9500
$ws.Range("M288").Value = 9267
$ws.Range("N288").Value = "$/caja 60 unidades"
$ws.Range("O288").Value = "Región de Arica y Parinacota"
$ws.Range("P288").Value = 154
$ws.Range("Q288").Value = 60
$ws.Range("R288").Value = "Hortaliza"
